$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.5586792620790276
$ws.Range("J2").Value = 0.5586792620790275
$ws.Range("M2").Value = 6.322177333333333
$ws.Range("N2").Value = 18.966532
$ws.Range("O2").Value = 0.08271011762055308
$ws.Range("P2").Value = 0.08271011762055309
$ws.Range("Q2").Value = 0.7803294891137778
$ws.Range("R2").Value = 7.022965402024
$ws.Range("S2").Value = 0.04620842747872018
$ws.Range("T2").Value = 0.04620842747872017

$ws.Range("I3").Value = 0.5586792620790276
$ws.Range("J3").Value = 0.5586792620790275
$ws.Range("O3").Value = 0.5401386314560596
$ws.Range("P3").Value = 0.5401386314560597
$ws.Range("S3").Value = 0.3017642520422472
$ws.Range("T3").Value = 0.3017642520422472

$ws.Range("I4").Value = 0.5586792620790276
$ws.Range("J4").Value = 0.5586792620790275
$ws.Range("M4").Value = 27.73243066666667
$ws.Range("N4").Value = 83.197292
$ws.Range("O4").Value = 0.3628105447549136
$ws.Range("P4").Value = 0.3628105447549136
$ws.Range("Q4").Value = 3.422939964038223
$ws.Range("R4").Value = 30.806459676344
$ws.Range("S4").Value = 0.2026947274181652
$ws.Range("T4").Value = 0.2026947274181651

$ws.Range("I5").Value = 0.5586792620790276
$ws.Range("J5").Value = 0.5586792620790275
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.096171666666667
$ws.Range("N5").Value = 3.288515
$ws.Range("O5").Value = 0.01434070616847367
$ws.Range("P5").Value = 0.01434070616847367
$ws.Range("Q5").Value = 0.1352975456922222
$ws.Range("R5").Value = 1.21767791123
$ws.Range("S5").Value = 0.008011855139895028
$ws.Range("T5").Value = 0.008011855139895026

$ws.Range("G6").Value = 0.09749966666666667
$ws.Range("H6").Value = 0.292499
$ws.Range("I6").Value = 0.4413207379209724
$ws.Range("J6").Value = 0.4413207379209724
$ws.Range("M6").Value = 6.322177333333333
$ws.Range("N6").Value = 18.966532
$ws.Range("O6").Value = 0.08271011762055308
$ws.Range("P6").Value = 0.08271011762055309
$ws.Range("Q6").Value = 0.6164101826075555
$ws.Range("R6").Value = 5.547691643468
$ws.Range("S6").Value = 0.03650169014183291
$ws.Range("T6").Value = 0.03650169014183291

$ws.Range("G7").Value = 0.09749966666666667
$ws.Range("H7").Value = 0.292499
$ws.Range("I7").Value = 0.4413207379209724
$ws.Range("J7").Value = 0.4413207379209724
$ws.Range("O7").Value = 0.5401386314560596
$ws.Range("P7").Value = 0.5401386314560597
$ws.Range("Q7").Value = 4.025468250168333
$ws.Range("R7").Value = 36.229214251515
$ws.Range("S7").Value = 0.2383743794138124
$ws.Range("T7").Value = 0.2383743794138124

$ws.Range("G8").Value = 0.09749966666666667
$ws.Range("H8").Value = 0.292499
$ws.Range("I8").Value = 0.4413207379209724
$ws.Range("J8").Value = 0.4413207379209724
$ws.Range("M8").Value = 27.73243066666667
$ws.Range("N8").Value = 83.197292
$ws.Range("O8").Value = 0.3628105447549136
$ws.Range("P8").Value = 0.3628105447549136
$ws.Range("Q8").Value = 2.703902745856445
$ws.Range("R8").Value = 24.335124712708
$ws.Range("S8").Value = 0.1601158173367485
$ws.Range("T8").Value = 0.1601158173367484

$ws.Range("G9").Value = 0.09749966666666667
$ws.Range("H9").Value = 0.292499
$ws.Range("I9").Value = 0.4413207379209724
$ws.Range("J9").Value = 0.4413207379209724
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.096171666666667
$ws.Range("N9").Value = 3.288515
$ws.Range("O9").Value = 0.01434070616847367
$ws.Range("P9").Value = 0.01434070616847367
$ws.Range("Q9").Value = 0.1068763721094445
$ws.Range("R9").Value = 0.9618873489850002
$ws.Range("S9").Value = 0.006328851028578639
$ws.Range("T9").Value = 0.006328851028578639
